# Daily automatic update of the electricity spot-price table (row 2)
# Source: Atualização automática de preços de eletricidade

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46027
$ws.Range("B2").Value = 88.18000000000001
$ws.Range("C2").Value = 71.66
$ws.Range("D2").Value = 62.99
$ws.Range("E2").Value = 61.93
$ws.Range("F2").Value = 61.93
$ws.Range("G2").Value = 64.20999999999999
$ws.Range("H2").Value = 77.28
$ws.Range("I2").Value = 94.12
$ws.Range("J2").Value = 102.86
$ws.Range("K2").Value = 107.92
$ws.Range("L2").Value = 104.34
$ws.Range("M2").Value = 100.52
$ws.Range("N2").Value = 97.70999999999999
$ws.Range("O2").Value = 98.04000000000001
$ws.Range("P2").Value = 97.98
$ws.Range("Q2").Value = 97.59999999999999
$ws.Range("R2").Value = 103.07
$ws.Range("S2").Value = 114.9
$ws.Range("T2").Value = 123.69
$ws.Range("U2").Value = 116.9
$ws.Range("V2").Value = 116.72
$ws.Range("W2").Value = 113.04
$ws.Range("X2").Value = 109.7
$ws.Range("Y2").Value = 103.24
$ws.Range("Z2").Value = 95.44
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 114.64
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 120.3
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 114.88
$ws.Range("AG2").Value = "0h-7h"
